# InstallTracker.xlsx update
# - Added 6 more days of "Actual" install-count data (rows 81-86 / idx 79-84
#   in the chart caches), fixing the last two existing days (81,82) and
#   filling in four brand-new days (83-86) whose Daily-delta / 7-day-average
#   formulas had previously been left blank.
# - Moved the frozen-pane scroll position / active selection to where the
#   author was last working (row ~68, cell D88).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Corrected "Actual" counts for the two most recent existing days.
# ---------------------------------------------------------------------
$ws.Range("C81").Value = 1435
$ws.Range("C82").Value = 1454

# ---------------------------------------------------------------------
# Four new days of data. Column C is typed/pasted in, and D/E (Daily
# delta and 7-day rolling Average) pick up the same formulas used by
# every other row in the table (D = Cn-Cn-1, E = (Cn-Cn-7)/7).
# ---------------------------------------------------------------------
$ws.Range("C83").Value = 1455
$ws.Range("C83").NumberFormat = "#,##0"
$ws.Range("D83").Formula = "=C83-C82"
$ws.Range("D83").NumberFormat = "0"
$ws.Range("E83").Formula = "=(C83-C76)/7"

$ws.Range("C84").Value = 1491
$ws.Range("C84").NumberFormat = "#,##0"
$ws.Range("D84").Formula = "=C84-C83"
$ws.Range("D84").NumberFormat = "0"
$ws.Range("E84").Formula = "=(C84-C77)/7"

$ws.Range("C85").Value = 1535
$ws.Range("C85").NumberFormat = "#,##0"
$ws.Range("D85").Formula = "=C85-C84"
$ws.Range("D85").NumberFormat = "0"
$ws.Range("E85").Formula = "=(C85-C78)/7"

$ws.Range("C86").Value = 1578
$ws.Range("C86").NumberFormat = "#,##0"
$ws.Range("D86").Formula = "=C86-C85"
$ws.Range("D86").NumberFormat = "0"
$ws.Range("E86").Formula = "=(C86-C79)/7"

# ---------------------------------------------------------------------
# View state: scroll the frozen pane down to the new rows and leave the
# selection on D88, matching where editing left off.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 68
$win.ScrollColumn = 1
$ws.Range("D88").Select() | Out-Null
